$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A89").Value = "2023-12-08 07:15:31"
$ws.Range("B89").Value = 0.0006000000000000001

$ws.Range("A90").Value = "2023-12-08 07:15:39"
$ws.Range("B90").Value = 0.0004
